$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 512
$ws1.Range("F3").Value = 3432
$ws1.Range("F4").Value = 93
$ws1.Range("F5").Value = 675

# Sheet "全部类型" (fourth sheet) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 512
$ws4.Range("F3").Value = 3432
$ws4.Range("F4").Value = 93
$ws4.Range("F5").Value = 675
